$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to keep edited cells
# free of any style index after forcing their values to be stored as text.
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "'37.517.71"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "'  +6.22%  "
$ws.Range("E2").Style = $plainStyle
$ws.Range("D3").Value = "'2.046.14"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "'  +3.46%  "
$ws.Range("E3").Style = $plainStyle
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = $plainStyle
$ws.Range("D5").Value = "'252.06"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "'  +5.20%  "
$ws.Range("E5").Style = $plainStyle
$ws.Range("E6").Value = "'  +2.44%  "
$ws.Range("E6").Style = $plainStyle
$ws.Range("D7").Value = "'66.59"
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = "'  +19.48%  "
$ws.Range("E7").Style = $plainStyle
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = $plainStyle
$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = "'  +6.32%  "
$ws.Range("E9").Style = $plainStyle
$ws.Range("D10").Value = "'59.42"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "'  +0.38%  "
$ws.Range("E10").Style = $plainStyle
$ws.Range("D11").Value = "'0.0756"
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = "'  +4.78%  "
$ws.Range("E11").Style = $plainStyle
$ws.Range("E12").Value = "'  +1.20%  "
$ws.Range("E12").Style = $plainStyle
$ws.Range("D13").Value = "'0.912"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "'  +2.92%  "
$ws.Range("E13").Style = $plainStyle
$ws.Range("D14").Value = "'15.15"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "'  +7.03%  "
$ws.Range("E14").Style = $plainStyle
$ws.Range("D15").Value = "'2.345.05"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "'  +3.37%  "
$ws.Range("E15").Style = $plainStyle
$ws.Range("E16").Value = "'  +8.05%  "
$ws.Range("E16").Style = $plainStyle
$ws.Range("D17").Value = "'20.82"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "'  +23.06%  "
$ws.Range("E17").Style = $plainStyle
$ws.Range("D18").Value = "'2.041.32"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "'  +3.42%  "
$ws.Range("E18").Style = $plainStyle
$ws.Range("D19").Value = "'37.407.96"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "'  +6.14%  "
$ws.Range("E19").Style = $plainStyle
$ws.Range("D20").Value = "'73.54"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "'  +5.58%  "
$ws.Range("E20").Style = $plainStyle
$ws.Range("D21").Value = "'0.0₃0878"
$ws.Range("D21").Style = $plainStyle
$ws.Range("E22").Value = "'  +7.69%  "
$ws.Range("E22").Style = $plainStyle
$ws.Range("D23").Value = "'237.52"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "'  +3.07%  "
$ws.Range("E23").Style = $plainStyle
$ws.Range("D24").Value = "'2.72"
$ws.Range("D24").Style = $plainStyle
$ws.Range("E25").Value = "'  -0.13%  "
$ws.Range("E25").Style = $plainStyle
$ws.Range("E26").Value = "'  +4.12%  "
$ws.Range("E26").Style = $plainStyle
$ws.Range("E27").Value = "'  +6.52%  "
$ws.Range("E27").Style = $plainStyle
$ws.Range("D28").Value = "'165.63"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "'  +1.80%  "
$ws.Range("E28").Style = $plainStyle
$ws.Range("D29").Value = "'19.91"
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = "'  +3.24%  "
$ws.Range("E29").Style = $plainStyle
$ws.Range("D30").Value = "'5.25"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "'  +11.23%  "
$ws.Range("E30").Style = $plainStyle
$ws.Range("D31").Value = "'0.122"
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = "'  +3.61%  "
$ws.Range("E31").Style = $plainStyle
$ws.Range("E32").Value = "'  +8.02%  "
$ws.Range("E32").Style = $plainStyle
$ws.Range("E33").Value = "'  +22.82%  "
$ws.Range("E33").Style = $plainStyle
$ws.Range("D34").Value = "'4.75"
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = "'  +12.55%  "
$ws.Range("E34").Style = $plainStyle
$ws.Range("E35").Value = "'  +5.79%  "
$ws.Range("E35").Style = $plainStyle
$ws.Range("E36").Value = "'  +10.24%  "
$ws.Range("E36").Style = $plainStyle
$ws.Range("D38").Value = "'6.05"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "'  +25.65%  "
$ws.Range("E38").Style = $plainStyle
$ws.Range("E39").Value = "'  +1.40%  "
$ws.Range("E39").Style = $plainStyle
$ws.Range("E40").Value = "'  +18.18%  "
$ws.Range("E40").Style = $plainStyle
$ws.Range("E42").Value = "'  +5.53%  "
$ws.Range("E42").Style = $plainStyle
$ws.Range("E43").Value = "'  +6.13%  "
$ws.Range("E43").Style = $plainStyle
$ws.Range("B44").Value = "'RenderToken"
$ws.Range("B44").Style = $plainStyle
$ws.Range("C44").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C44").Style = $plainStyle
$ws.Range("D44").Value = "'2.73"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "'  +21.96%  "
$ws.Range("E44").Style = $plainStyle
$ws.Range("B45").Value = "'ARBITRUM"
$ws.Range("B45").Style = $plainStyle
$ws.Range("C45").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C45").Style = $plainStyle
$ws.Range("D45").Value = "'1.14"
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = "'  +6.74%  "
$ws.Range("E45").Style = $plainStyle
$ws.Range("D46").Value = "'8.11"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "'  +10.06%  "
$ws.Range("E46").Style = $plainStyle
$ws.Range("E47").Value = "'  +11.54%  "
$ws.Range("E47").Style = $plainStyle
$ws.Range("D48").Value = "'95.32"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "'  +6.63%  "
$ws.Range("E48").Style = $plainStyle
$ws.Range("D49").Value = "'1.427.31"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "'  +5.82%  "
$ws.Range("E49").Style = $plainStyle
$ws.Range("E50").Value = "'  +2.70%  "
$ws.Range("E50").Style = $plainStyle
$ws.Range("D51").Value = "'47.55"
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = "'  +5.45%  "
$ws.Range("E51").Style = $plainStyle
